# "cleaned up state comma spaces"
# Remove the space after the comma in the "states" column values on Sheet1,
# and make Sheet1 (with cell B7 selected) the active sheet/selection instead
# of Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("B2").Value = "Wyoming,Michigan"
$ws1.Range("B3").Value = "Wisconsin,Nevada,California"
$ws1.Range("B4").Value = "Florida,Washington"
$ws1.Range("B6").Value = "Washington,Oregon,California"

# Move the active sheet/selection to Sheet1 (cell B7), matching the
# workbook's saved view state after the edit.
$ws1.Activate()
$ws1.Range("B7").Select() | Out-Null
